$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at position 230 (shifts existing rows 230-246 down to 234-250)
$ws.Range("A230:A233").EntireRow.Insert()

# Row 230
$ws.Range("A230").Value = 2
$ws.Range("B230").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C230").Value = "Coquimbo"
$ws.Range("D230").Value = 44783
$ws.Range("E230").Value = 4
$ws.Range("F230").Value = 100112013
$ws.Range("G230").Value = "Alcachofa"
$ws.Range("H230").Value = "Argentina(o)"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 1100
$ws.Range("K230").Value = 10000
$ws.Range("L230").Value = 11000
$ws.Range("M230").Value = 10500
$ws.Range("N230").Value = "`$/caja 50 unidades"
$ws.Range("O230").Value = "Provincia de Limarí"
$ws.Range("P230").Value = 210
$ws.Range("Q230").Value = 50
$ws.Range("R230").Value = "Hortaliza"

# Row 231
$ws.Range("A231").Value = 2
$ws.Range("B231").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C231").Value = "Coquimbo"
$ws.Range("D231").Value = 44783
$ws.Range("E231").Value = 4
$ws.Range("F231").Value = 100112013
$ws.Range("G231").Value = "Alcachofa"
$ws.Range("H231").Value = "Argentina(o)"
$ws.Range("I231").Value = "Segunda"
$ws.Range("J231").Value = 700
$ws.Range("K231").Value = 8000
$ws.Range("L231").Value = 9000
$ws.Range("M231").Value = 8500
$ws.Range("N231").Value = "`$/caja 70 unidades"
$ws.Range("O231").Value = "Provincia de Limarí"
$ws.Range("P231").Value = 121
$ws.Range("Q231").Value = 70
$ws.Range("R231").Value = "Hortaliza"

# Row 232
$ws.Range("A232").Value = 2
$ws.Range("B232").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C232").Value = "Coquimbo"
$ws.Range("D232").Value = 44783
$ws.Range("E232").Value = 4
$ws.Range("F232").Value = 100112013
$ws.Range("G232").Value = "Alcachofa"
$ws.Range("H232").Value = "Española"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 1900
$ws.Range("K232").Value = 11000
$ws.Range("L232").Value = 12000
$ws.Range("M232").Value = 11500
$ws.Range("N232").Value = "`$/caja 30 unidades"
$ws.Range("O232").Value = "Provincia de Limarí"
$ws.Range("P232").Value = 383
$ws.Range("Q232").Value = 30
$ws.Range("R232").Value = "Hortaliza"

# Row 233
$ws.Range("A233").Value = 2
$ws.Range("B233").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C233").Value = "Coquimbo"
$ws.Range("D233").Value = 44783
$ws.Range("E233").Value = 4
$ws.Range("F233").Value = 100112013
$ws.Range("G233").Value = "Alcachofa"
$ws.Range("H233").Value = "Madrigal"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 300
$ws.Range("K233").Value = 9000
$ws.Range("L233").Value = 10000
$ws.Range("M233").Value = 9500
$ws.Range("N233").Value = "`$/caja 40 unidades"
$ws.Range("O233").Value = "Provincia de Limarí"
$ws.Range("P233").Value = 238
$ws.Range("Q233").Value = 40
$ws.Range("R233").Value = "Hortaliza"
